$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "ByCityName"  (grows from 3 rows to 9 rows)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ByCityName")

$ws1.Range("B2").Value = "Existing city. No country code. No Other features"
$ws1.Range("B3").Value = "Existing city. Valid country code. No Other features"
$ws1.Range("I3").Value = "DEFAULT"

$ws1.Range("A4").Formula = "=A3+1"
$ws1.Range("B4").Value = "Existing city. Valid country code. Language english"
$ws1.Range("C4").Value = "CORRECT"
$ws1.Range("D4").Value = "London"
$ws1.Range("E4").Value = "UK"
$ws1.Range("F4").Value = "JSON"
$ws1.Range("G4").Value = "DEFAULT"
$ws1.Range("H4").Value = "STANDARD"
$ws1.Range("I4").Value = "ENGLISH"
$ws1.Range("J4").Value = 200
$ws1.Range("L4").Value = 2643743
$ws1.Range("M4").Value = "London"

$ws1.Range("A5").Formula = "=A4+1"
$ws1.Range("B5").Value = "Invalid city. Invalid country code. Search Acurate. Language german"
$ws1.Range("C5").Value = "CORRECT"
$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 99
$ws1.Range("F5").Value = "JSON"
$ws1.Range("G5").Value = "ACCURATE"
$ws1.Range("H5").Value = "STANDARD"
$ws1.Range("I5").Value = "GERMAN"
$ws1.Range("J5").Value = 200
$ws1.Range("L5").Value = 8182072
$ws1.Range("M5").Value = "Estanzuela"

$ws1.Range("A6").Formula = "=A5+1"
$ws1.Range("B6").Value = "Existing city. Valid country code. Unit Celsius. Language english"
$ws1.Range("C6").Value = "CORRECT"
$ws1.Range("D6").Value = "London"
$ws1.Range("E6").Value = "UK"
$ws1.Range("F6").Value = "JSON"
$ws1.Range("G6").Value = "DEFAULT"
$ws1.Range("H6").Value = "METRIC"
$ws1.Range("I6").Value = "ENGLISH"
$ws1.Range("J6").Value = 200
$ws1.Range("L6").Value = 2643743
$ws1.Range("M6").Value = "London"

$ws1.Range("A7").Formula = "=A6+1"
$ws1.Range("B7").Value = "Existing city. Valid country code. Unit Fahrenheit. Language english"
$ws1.Range("C7").Value = "CORRECT"
$ws1.Range("D7").Value = "London"
$ws1.Range("E7").Value = "UK"
$ws1.Range("F7").Value = "JSON"
$ws1.Range("G7").Value = "DEFAULT"
$ws1.Range("H7").Value = "IMPERIAL"
$ws1.Range("I7").Value = "ENGLISH"
$ws1.Range("J7").Value = 200
$ws1.Range("L7").Value = 2643743
$ws1.Range("M7").Value = "London"

$ws1.Range("A8").Formula = "=A7+1"
$ws1.Range("B8").Value = "Existing city. Valid country code. Search like, Unit Fahrenheit. Language english"
$ws1.Range("C8").Value = "CORRECT"
$ws1.Range("D8").Value = "London"
$ws1.Range("E8").Value = "UK"
$ws1.Range("F8").Value = "JSON"
$ws1.Range("G8").Value = "LIKE"
$ws1.Range("H8").Value = "IMPERIAL"
$ws1.Range("I8").Value = "ENGLISH"
$ws1.Range("J8").Value = 200
$ws1.Range("L8").Value = 2643743
$ws1.Range("M8").Value = "London"

$ws1.Range("A9").Formula = "=A8+1"
$ws1.Range("B9").Value = "City Part. Search like, Unit Celsius. Language english"
$ws1.Range("C9").Value = "CORRECT"
$ws1.Range("D9").Value = "Mex"
$ws1.Range("F9").Value = "JSON"
$ws1.Range("G9").Value = "LIKE"
$ws1.Range("H9").Value = "METRIC"
$ws1.Range("I9").Value = "ENGLISH"
$ws1.Range("J9").Value = 200
$ws1.Range("L9").Value = 3530597
$ws1.Range("M9").Value = "Mexico City"

# column D needs a bit more width on this sheet now that longer city
# names/zip parts live there (mirrors the author's column width tweak)
$ws1.Columns.Item(4).ColumnWidth = 14.28515625

# ---------------------------------------------------------------
# Sheet "ById"  (grows from 3 rows to 5 rows)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ById")

$ws2.Range("B2").Value = "Valid Id. No other Features."
$ws2.Range("B3").Value = "Invalid Id"

$ws2.Range("A4").Formula = "=A3+1"
$ws2.Range("B4").Value = "Valid Id. Unit Celsius"
$ws2.Range("C4").Value = "CORRECT"
$ws2.Range("D4").Value = 8182072
$ws2.Range("E4").Value = "JSON"
$ws2.Range("F4").Value = "DEFAULT"
$ws2.Range("G4").Value = "METRIC"
$ws2.Range("H4").Value = "DEFAULT"
$ws2.Range("I4").Value = 200
$ws2.Range("K4").Value = 8182072
$ws2.Range("L4").Value = "Estanzuela"

$ws2.Range("A5").Formula = "=A4+1"
$ws2.Range("B5").Value = "Valid Id. Unit Celsius. Language Spanish"
$ws2.Range("C5").Value = "CORRECT"
$ws2.Range("D5").Value = 3996063
$ws2.Range("E5").Value = "JSON"
$ws2.Range("F5").Value = "DEFAULT"
$ws2.Range("G5").Value = "METRIC"
$ws2.Range("H5").Value = "SPANISH"
$ws2.Range("I5").Value = 200
$ws2.Range("K5").Value = 3996063
$ws2.Range("L5").Value = "Mexico"

# Copy the "Expected Name" cell style (vertical-center) down onto the two
# newly-added rows, same style used by the existing rows in this column.
$ws2.Range("L2").Copy()
$ws2.Range("L4").PasteSpecial(-4122)
$ws2.Range("L2").Copy()
$ws2.Range("L5").PasteSpecial(-4122)

$ws2.Range("L2").Copy()
$ws1.Range("M9").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet "ByZipCode"  (same two rows, the TestCaseId numbers move)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ByZipCode")
$ws4.Range("A2").Value = 4001

# ---------------------------------------------------------------
# Selections / active tab - move the cursor the same way the author
# left it in each sheet, then bring "ByCityName" to the front tab.
# ---------------------------------------------------------------
$ws2.Range("K5:L5").Select()

$ws3 = $wb.Worksheets.Item("ByGeoCoords")
$ws3.Range("I8").Select()

$ws4.Range("B3").Select()

$ws1.Activate()
